$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Key/Value rows (27-33) appended after the existing "files" row (26),
# for the network log viewer feature (initial network log display).
$rows = @(
    @("network_log_title", "Network Log Viewer"),
    @("network_active_connections", "Active Connections"),
    @("network_log_proto", "Proto"),
    @("network_log_local_addy", "Local Address"),
    @("network_log_foreign_addy", "Foreign Address"),
    @("network_log_state", "State"),
    @("network_log_pid", "PID")
)

$startRow = 27
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Update the view: scroll/select to reflect the newly added rows at the
# bottom of the sheet (selection moves from A26 to A34).
$lastRow = $startRow + $rows.Length
$ws.Range("A" + $lastRow).Select()
